# B6-PowerPoint.pptx edit — Mon, Jul 20, 2020  8:07:15 AM
#
# 1) Re-style the three tables (slides 14, 15, 16) from the bespoke
#    "Table_0" style to the built-in table style {8ACB93C8-694C-4A17-
#    801B-8429A476A4FE}.
# 2) Re-colour the deck's theme (slide master) from the "Integral" /
#    "Red Violet" palette to the standard "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Tables -------------------------------------------------------
$newTableStyle = "{8ACB93C8-694C-4A17-801B-8429A476A4FE}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2. Theme colours --------------------------------------------------
# PowerPoint ThemeColorScheme.Colors() index order:
#  1 dk1, 2 lt1, 3 dk2, 4 lt2,
#  5 accent1, 6 accent2, 7 accent3, 8 accent4, 9 accent5, 10 accent6,
#  11 hlink, 12 folHlink
function Set-ThemeColor($scheme, $index, $r, $g, $b) {
    $scheme.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

Set-ThemeColor $colorScheme 1  0x00 0x00 0x00   # dk1
Set-ThemeColor $colorScheme 2  0xFF 0xFF 0xFF   # lt1
Set-ThemeColor $colorScheme 3  0x44 0x54 0x6A   # dk2
Set-ThemeColor $colorScheme 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeColor $colorScheme 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeColor $colorScheme 6  0xED 0x7D 0x31   # accent2
Set-ThemeColor $colorScheme 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeColor $colorScheme 8  0xFF 0xC0 0x00   # accent4
Set-ThemeColor $colorScheme 9  0x44 0x72 0xC4   # accent5
Set-ThemeColor $colorScheme 10 0x70 0xAD 0x47   # accent6
Set-ThemeColor $colorScheme 11 0x05 0x63 0xC1   # hlink
Set-ThemeColor $colorScheme 12 0x95 0x4F 0x72   # folHlink
